$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'89.092.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.12%  '

$ws.Range("D3").Value = "'3.136.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.05%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = "'215.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.92%  '

$ws.Range("D6").Value = "'635.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.15%  '

$ws.Range("D7").Value = "'0.397"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.98%  '

$ws.Range("D8").Value = "'0.765"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.27%  '

$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").Value = "'3.133.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.18%  '

$ws.Range("E11").Value = '  -5.03%  '

$ws.Range("E12").Value = '  -0.50%  '

$ws.Range("E13").Value = '  -5.50%  '

$ws.Range("E14").Value = '  -0.40%  '

$ws.Range("D15").Value = "'88.880.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.12%  '

$ws.Range("D16").Value = "'3.712.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.19%  '

$ws.Range("D17").Value = "'32.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.32%  '

$ws.Range("D18").Value = "'3.145.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.16%  '

$ws.Range("E19").Value = '  +18.15%  '

$ws.Range("D20").Value = "'3.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.20%  '

$ws.Range("D21").Value = "'13.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.09%  '

$ws.Range("D22").Value = "'428.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.36%  '

$ws.Range("D23").Value = "'8.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.64%  '

$ws.Range("D24").Value = "'4.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.50%  '

$ws.Range("D25").Value = "'5.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.60%  '

$ws.Range("D26").Value = "'83.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.76%  '

$ws.Range("D27").Value = "'11.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.27%  '

$ws.Range("D28").Value = "'3.301.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.38%  '

$ws.Range("E30").Value = '  -12.71%  '

$ws.Range("D31").Value = "'0.995"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.56%  '

$ws.Range("D32").Value = "'4.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +10.36%  '

$ws.Range("E33").Value = '  -6.32%  '

$ws.Range("D34").Value = "'509.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.11%  '

$ws.Range("D35").Value = "'0.148"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +13.54%  '

$ws.Range("D36").Value = "'7.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.14%  '

$ws.Range("D37").Value = "'1.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.23%  '

$ws.Range("E38").Value = '  -4.51%  '

$ws.Range("D39").Value = "'21.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.83%  '

$ws.Range("E40").Value = '  -0.75%  '

$ws.Range("E41").Value = '  +0.19%  '

$ws.Range("E42").Value = '  -0.11%  '

$ws.Range("D43").Value = "'1.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.52%  '

$ws.Range("E44").Value = '  -7.39%  '

$ws.Range("D45").Value = "'145.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.84%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = "'0.131"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.17%  '

$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").Value = "'43.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.98%  '

$ws.Range("D48").Value = "'164.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.29%  '

$ws.Range("E49").Value = '  -1.35%  '

$ws.Range("D50").Value = "'24.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.03%  '

$ws.Range("B51").Value = 'ImmutableX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D51").Value = "'1.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.80%  '
